# Updates the cryptos price list (columns D = Price, E = Volume(1h))
# for the data rows whose figures changed in this refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.900.03"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.765.38"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.Value = "'329.04"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("E6").Value = "  +0.01%  "
$c = $ws.Range("D7")
$c.Value = "'0.4539"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("E8").Value = "  -1.57%  "
$c = $ws.Range("D9")
$c.Value = "'42.03"
$c.ClearFormats()
$ws.Range("E9").Value = "  +1.40%  "
$c = $ws.Range("D10")
$c.Value = "'0.07386"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.13%  "
$c = $ws.Range("D11")
$c.Value = "'1.097"
$c.ClearFormats()
$ws.Range("E11").Value = "  +1.11%  "
$c = $ws.Range("D12")
$c.Value = "'1.002"
$c.ClearFormats()
$ws.Range("E12").Value = "  +0.07%  "
$c = $ws.Range("D13")
$c.Value = "'20.73"
$c.ClearFormats()
$ws.Range("E13").Value = "  -0.24%  "
$c = $ws.Range("D14")
$c.Value = "'6.004"
$c.ClearFormats()
$ws.Range("E14").Value = "  +0.23%  "
$c = $ws.Range("D15")
$c.Value = "'7.192"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "1.765.61"
$ws.Range("E16").Value = "  +0.53%  "
$c = $ws.Range("D17")
$c.Value = "'92.35"
$c.ClearFormats()
$ws.Range("E17").Value = "  -2.27%  "
$c = $ws.Range("D18")
$c.Value = "'0.00001057"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.13%  "
$c = $ws.Range("D19")
$c.Value = "'0.06445"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("E20").Value = "  +0.06%  "
$c = $ws.Range("D21")
$c.Value = "'17.00"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "27.926.29"
$ws.Range("E23").Value = "  +1.13%  "
$c = $ws.Range("D24")
$c.Value = "'11.25"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.55%  "
$c = $ws.Range("D25")
$c.Value = "'2.153"
$c.ClearFormats()
$ws.Range("E25").Value = "  +3.44%  "
$c = $ws.Range("D26")
$c.Value = "'161.80"
$c.ClearFormats()
$ws.Range("E26").Value = "  -2.35%  "
$c = $ws.Range("D27")
$c.Value = "'20.14"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "1.971.35"
$ws.Range("E28").Value = "  +0.80%  "
$c = $ws.Range("D29")
$c.Value = "'2.172"
$c.ClearFormats()
$ws.Range("E29").Value = "  +2.31%  "
$c = $ws.Range("D30")
$c.Value = "'124.01"
$c.ClearFormats()
$ws.Range("E30").Value = "  -1.17%  "
$c = $ws.Range("D31")
$c.Value = "'1.078"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.87%  "
$c = $ws.Range("D32")
$c.Value = "'0.09274"
$c.ClearFormats()
$ws.Range("E32").Value = "  +1.18%  "
$c = $ws.Range("D33")
$c.Value = "'5.590"
$c.ClearFormats()
$ws.Range("E33").Value = "  +1.51%  "
$c = $ws.Range("D34")
$c.Value = "'3.647"
$c.ClearFormats()
$ws.Range("E34").Value = "  -0.08%  "
$c = $ws.Range("D35")
$c.Value = "'11.82"
$c.ClearFormats()
$ws.Range("E35").Value = "  +0.95%  "
$c = $ws.Range("D36")
$c.Value = "'0.02272"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.51%  "
$c = $ws.Range("D37")
$c.Value = "'0.06125"
$c.ClearFormats()
$ws.Range("E37").Value = "  +2.10%  "
$c = $ws.Range("D38")
$c.Value = "'0.2088"
$c.ClearFormats()
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -0.81%  "
$c = $ws.Range("D43")
$c.Value = "'7.839"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.63%  "
$c = $ws.Range("D44")
$c.Value = "'13.19"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.29%  "
$c = $ws.Range("D45")
$c.Value = "'3.735"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.51%  "
$c = $ws.Range("D46")
$c.Value = "'0.5842"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.22%  "
$c = $ws.Range("D47")
$c.Value = "'122.83"
$c.ClearFormats()
$ws.Range("E47").Value = "  +0.82%  "
$c = $ws.Range("D48")
$c.Value = "'1.937"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.20%  "
$c = $ws.Range("D49")
$c.Value = "'1.130"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.06%  "
$c = $ws.Range("D50")
$c.Value = "'0.06809"
$c.ClearFormats()
$ws.Range("E50").Value = "  -1.12%  "
$c = $ws.Range("D51")
$c.Value = "'72.94"
$c.ClearFormats()
$ws.Range("E51").Value = "  +2.06%  "
